$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows right before row 961 (shifts old rows 961..1068 down to 963..1070).
$ws.Range("A961:A962").EntireRow.Insert()

# Seed the two new rows from the rows that just landed below them (963/964 carry the
# same fixed columns: A,B,C,E,F,G,H,I,N,Q,R) so formatting/shared columns match exactly.
$ws.Rows.Item(963).Copy()
$ws.Rows.Item(961).PasteSpecial()
$ws.Rows.Item(964).Copy()
$ws.Rows.Item(962).PasteSpecial()

# Now overwrite the columns that actually hold new data for this new weekly entry.
$ws.Range("D961").Value = 45212
$ws.Range("J961").Value = 1600
$ws.Range("K961").Value = 800
$ws.Range("L961").Value = 900
$ws.Range("M961").Value = 850
$ws.Range("O961").Value = "Región Metropolitana"
$ws.Range("P961").Value = 850

$ws.Range("D962").Value = 45212
$ws.Range("J962").Value = 970
$ws.Range("K962").Value = 700
$ws.Range("L962").Value = 700
$ws.Range("M962").Value = 700
$ws.Range("O962").Value = "Región Metropolitana"
$ws.Range("P962").Value = 700
